# Change table header keyword from id='...' to class='...' across all
# ObjTables "Data" table definition sheets, so custom ids can later be used
# to track the table origin of each model instance.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets()) {
    foreach ($cell in $ws.UsedRange().Cells()) {
        $v = $cell.Value()
        if ($v -ne $null -and $v -like "!!ObjTables type=*id=*") {
            $nv = $v -replace "' id='", "' class='"
            $cell.Value = $nv
        }
    }
}
